$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B, C, D, E for rows 2-51 (coin name, link, price, volume)
$data = @(
    @("Bitcoin", "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc", "25.856.17", "  -0.19%  "),
    @("Ethereum", "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth", "1.738.95", "  -0.15%  "),
    @("TetherUSD", "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt", "0.9995", "  +0.01%  "),
    @("BNB", "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb", "241.13", "  +4.38%  "),
    @("USDC", "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc", "0.9996", "  +0.03%  "),
    @("XRP", "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp", "0.5222", "  -0.76%  "),
    @("Cardano", "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada", "0.2746", "  -0.88%  "),
    @("OKB", "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb", "39.94", "  +1.12%  "),
    @("Dogecoin", "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge", "0.06164", "  +0.17%  "),
    @("WrappedEther", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth", "1.740.70", "  +0.03%  "),
    @("TRON", "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx", "0.07180", "  +0.90%  "),
    @("Solana", "https://coinranking.com/coin/zNZHO_Sjf+solana-sol", "15.00", "  -1.57%  "),
    @("Polygon", "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic", "0.6421", "  -0.70%  "),
    @("Polkadot", "https://coinranking.com/coin/25W7FG7om+polkadot-dot", "4.615", "  +1.86%  "),
    @("Litecoin", "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc", "77.40", "  +0.42%  "),
    @("Dai", "https://coinranking.com/coin/MoTuySvg7+dai-dai", "1.000", "  +0.07%  "),
    @("BinanceUSD", "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd", "0.9997", "  +0.09%  "),
    @("WrappedBTC", "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc", "25.875.72", "  +0.01%  "),
    @("Avalanche", "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax", "11.72", "  +1.54%  "),
    @("ShibaInu", "https://coinranking.com/coin/xz24e0BjL+shibainu-shib", "0.000006774", "  +1.42%  "),
    @("WrappedliquidstakedEther2.0", "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth", "1.961.62", "  +0.11%  "),
    @("Uniswap", "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni", "4.270", "  -0.12%  "),
    @("Cosmos", "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom", "8.641", "  -1.86%  "),
    @("Chainlink", "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link", "5.266", "  +1.73%  "),
    @("Monero", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr", "138.34", "  -1.47%  "),
    @("Toncoin", "https://coinranking.com/coin/67YlI0K1b+toncoin-ton", "1.518", "  -0.13%  "),
    @("EthereumClassic", "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc", "15.19", "  -0.08%  "),
    @("LidoDAOToken", "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo", "1.767", "  -1.85%  "),
    @("BitcoinCash", "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch", "105.28", "  +2.81%  "),
    @("InternetComputer(DFINITY)", "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp", "3.956", "  +5.76%  "),
    @("Stellar", "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm", "0.08276", "  -0.72%  "),
    @("Filecoin", "https://coinranking.com/coin/ymQub4fuB+filecoin-fil", "3.679", "  +2.46%  "),
    @("Hedera", "https://coinranking.com/coin/jad286TjB+hedera-hbar", "0.04625", "  +2.51%  "),
    @("HuobiToken", "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht", "2.642", "  +1.18%  "),
    @("ARBITRUM", "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb", "0.9874", "  +1.17%  "),
    @("ImmutableX", "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx", "0.6172", "  -0.58%  "),
    @("MXToken", "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx", "2.683", "  -0.11%  "),
    @("VeChain", "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet", "0.01602", "  +1.17%  "),
    @("RenderToken", "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr", "1.927", "  +0.86%  "),
    @("PaxDollar", "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp", "0.9994", "  +0.05%  "),
    @("Quant", "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt", "98.86", "  -1.12%  "),
    @("TheSandbox", "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand", "0.3847", "  -0.55%  "),
    @("TrustWalletToken", "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt", "0.7432", "  +1.81%  "),
    @("FraxShare", "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs", "4.992", "  -0.49%  "),
    @("Algorand", "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo", "0.1128", "  +0.63%  "),
    @("Aptos", "https://coinranking.com/coin/HGYj5JCv5+aptos-apt", "6.246", "  -0.01%  "),
    @("Cronos", "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro", "0.05243", "  -1.58%  "),
    @("Aave", "https://coinranking.com/coin/ixgUfzmLR+aave-aave", "54.69", "  +1.93%  "),
    @("Elrond", "https://coinranking.com/coin/omwkOTglq+elrond-egld", "30.45", "  +1.06%  "),
    @("EnergySwap", "https://coinranking.com/coin/SbWqqTui-+energyswap-ens", "7.594", "  -1.00%  ")
)

$row = 2
foreach ($item in $data) {
    $cB = $ws.Cells.Item($row, 2)
    $cB.Value = $item[0]

    $cC = $ws.Cells.Item($row, 3)
    $cC.Value = $item[1]

    # D and E hold numeric-looking text (e.g. "0.9995", "25.856.17") that must
    # stay stored as text, matching the source file's inlineStr cells.
    $cD = $ws.Cells.Item($row, 4)
    $cD.NumberFormat = "@"
    $cD.Value = $item[2]
    $cD.NumberFormat = "General"
    $cD.Style = "Normal"

    $cE = $ws.Cells.Item($row, 5)
    $cE.NumberFormat = "@"
    $cE.Value = $item[3]
    $cE.NumberFormat = "General"
    $cE.Style = "Normal"

    $row = $row + 1
}
